$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Bring this sheet into focus (it becomes the active/selected tab)
$ws.Activate()

# Insert a new (blank) column before the existing "Late" column (column N),
# shifting Late / heading / Outstanding one column to the right (N->O, O->P, P->Q)
$ws.Columns.Item(14).Insert()

# Excel carries the width of the column to the left (M) onto the freshly
# inserted column when doing a column insert
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

# Final selection left on the sheet after the edit
$ws.Range("K17").Select() | Out-Null
